$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 worker/period rows that are no longer part of the account
# statement (this also shifts the last table row, with its bottom-border
# styling, and the footer rows up into place).
$ws.Rows("19:23").Delete()

# Refresh the remaining 4 data rows with the updated worker information.
$ws.Range("C16").Value = "73162788"
$ws.Range("D16").Value = "MAYER ENRIQUE CANTILLO ALTAMIRANDA"
$ws.Range("E16").Value = "1806"
$ws.Range("F16").Value = 67776
$ws.Range("G16").Value = 781242

$ws.Range("C17").Value = "9095273"
$ws.Range("D17").Value = "JORGE ELIECER ROQUE JIMENEZ"
$ws.Range("E17").Value = "1806"
$ws.Range("F17").Value = 50832
$ws.Range("G17").Value = 781242

$ws.Range("C18").Value = "73146705"
$ws.Range("D18").Value = "BLAS ALBERTO TORRES MAZA"
$ws.Range("E18").Value = "1806"
$ws.Range("F18").Value = 50832
$ws.Range("G18").Value = 1200000

$ws.Range("C19").Value = "73192056"
$ws.Range("D19").Value = "JESUS MARIA CANTILLO ALTAMIRANDA"
$ws.Range("E19").Value = "1806"
$ws.Range("F19").Value = 67776
$ws.Range("G19").Value = 828116

# Update the summary figures at the top of the statement.
$ws.Range("E11").Value = 237216
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 1
